# ---------------------------------------------------------------------------
# Applies the "Completed Word doc report" edit to Report.docx
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Abstract paragraph: drop the stray "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2) Running Instructions > Phase 2 bullet (currently just "A").
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(11)
$p2Start = $p2.Range.Start
$p2.Range.Text = "To run Phase 2's code, simply call the python file named phase2.py. This will create intermediate files and the index files, being re.idx, te.idx, em.idx,  and da.idx. Example:"

$rng = $d.Range($p2Start, $d.Paragraphs(11).Range.End)
$rng.Find.Execute("phase2.py") | Out-Null
$rng.Italic = 1

$rng = $d.Range($p2Start, $d.Paragraphs(11).Range.End)
$rng.Find.Execute("re.idx") | Out-Null
$rng.Italic = 1

$rng = $d.Range($p2Start, $d.Paragraphs(11).Range.End)
$rng.Find.Execute(" te.idx") | Out-Null
$rng.Italic = 1

$rng = $d.Range($p2Start, $d.Paragraphs(11).Range.End)
$rng.Find.Execute(" em.idx") | Out-Null
$rng.Italic = 1

$rng = $d.Range($p2Start, $d.Paragraphs(11).Range.End)
$rng.Find.Execute(" da.idx") | Out-Null
$rng.Italic = 1

# New centred "Example:" line right after the bullet.
$d.Paragraphs(11).Range.InsertParagraphAfter()
$ex2 = $d.Paragraphs(12)
$ex2.Range.ListFormat.RemoveNumbers()
$ex2.Style = "Normal"
$ex2.Alignment = 1
$ex2.Range.Text = ">> python phase2.py"
$ex2.Range.Font.Name = "Consolas"
$ex2.Range.Font.Size = 12

# ---------------------------------------------------------------------------
# 3) Running Instructions > Phase 3 bullet (currently just "B").
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(13)
$p3Start = $p3.Range.Start
$p3.Range.Text = "To run Phase 3's code, and bring up the query interface, simply call the python file named phase3.py. Follow the instructions on screen or write queries given the query language. Example:"

$rng = $d.Range($p3Start, $d.Paragraphs(13).Range.End)
$rng.Find.Execute("phase3") | Out-Null
$rng.Italic = 1

# New centred "Example:" line right after the bullet.
$d.Paragraphs(13).Range.InsertParagraphAfter()
$ex3 = $d.Paragraphs(14)
$ex3.Range.ListFormat.RemoveNumbers()
$ex3.Style = "Normal"
$ex3.Alignment = 1
$ex3.Range.Text = ">> python phase3.py"
$ex3.Range.Font.Name = "Consolas"
$ex3.Range.Font.Size = 12

# ---------------------------------------------------------------------------
# 4) Testing Strategy > Phase 2 bullet (currently just "A").
# ---------------------------------------------------------------------------
$t2 = $d.Paragraphs(18)
$t2.SpaceBefore = 12
$t2Start = $t2.Range.Start
$t2.Range.Text = "Phase 2 was tested by looking for sorted files and indexed files given the files had keys and data separated. Using db_dump in the terminal on the index file to see the indices, we were able to verify it was done properly."

$rng = $d.Range($t2Start, $d.Paragraphs(18).Range.End)
$rng.Find.Execute("db_dump") | Out-Null
$rng.Italic = 1

# ---------------------------------------------------------------------------
# 5) Testing Strategy > Phase 3 bullet (currently just "B").
# ---------------------------------------------------------------------------
$d.Paragraphs(19).Range.Text = "To test Phase 3, the interface created was used, and a database was created according to the data retrieved specification. Queries were performed on data sets given."

# ---------------------------------------------------------------------------
# 6) Group Work Break-Down Strategy paragraph.
# ---------------------------------------------------------------------------
$gw = $d.Paragraphs(22)
$gwStart = $gw.Range.Start
$gwEnd = $d.Paragraphs(22).Range.End

$rng = $d.Range($gwStart, $gwEnd)
$rng.Find.Execute("approximately 3 hours", $true, $false, $false, $false, $false, $true, 1, $false, "approximately 5 hours", 2) | Out-Null

$rng = $d.Range($gwStart, $d.Paragraphs(22).Range.End)
$rng.Find.Execute("HERE and HERE hours, respectively", $true, $false, $false, $false, $false, $true, 1, $false, "3 and 10 hours, respectively", 2) | Out-Null

$rng = $d.Range($gwStart, $d.Paragraphs(22).Range.End)
$rng.Find.Execute("by both members of the team.", $true, $false, $false, $false, $false, $true, 1, $false, "by both members of the team in the final days of development.", 2) | Out-Null

# Strip any stray bold left over from the old "HERE" runs and restore normal weight.
$rng = $d.Range($gwStart, $d.Paragraphs(22).Range.End)
$rng.Font.Bold = 0

# Re-anchor the "_GoBack" bookmark at the point the edit ended (right before the
# final "Coordination was kept..." sentence).
$rng = $d.Range($gwStart, $d.Paragraphs(22).Range.End)
$rng.Find.Execute("by both members of the team in the final days of development") | Out-Null
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)
